$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# ------------------------------------------------------------------
# Capture values that need to be re-used verbatim (preserves exact
# unicode text / shared-string de-duplication) before we start
# shuffling rows around.
# ------------------------------------------------------------------
$demandNLText = $ws.Range("C18").Value2

# ------------------------------------------------------------------
# 1) Split "fix_demand_and_profiles_to_initial_year" (row 18) into two
#    separate rows: "fix_demand_to_initial_year" (row 18, FALSE) and a
#    brand new "fix_profiles_to_initial_year" (row 19, TRUE). This is
#    done by inserting one new row before the old row 19, which pushes
#    the old rows 19-25 down to 20-26.
# ------------------------------------------------------------------
$ws.Rows.Item(19).Insert()

# ------------------------------------------------------------------
# 2) Update B17 (fix_price_year) from 2020 to 2030
# ------------------------------------------------------------------
$ws.Range("B17").Value2 = 2030

# ------------------------------------------------------------------
# 3) Rewrite row 18 as "fix_demand_to_initial_year" = FALSE
#    (C18 keeps its original text/value, no change needed there)
# ------------------------------------------------------------------
$ws.Range("A18").Value2 = "fix_demand_to_initial_year"
$ws.Range("B18").Value2 = $false

# ------------------------------------------------------------------
# 4) Populate the newly inserted row 19 as
#    "fix_profiles_to_initial_year" = TRUE, reusing the same note text
#    that row 18 originally had.
# ------------------------------------------------------------------
$ws.Range("A19").Value2 = "fix_profiles_to_initial_year"
$ws.Range("B19").Value2 = $true
$ws.Range("C19").Value2 = $demandNLText

# ------------------------------------------------------------------
# 5) Rows that used to be 19-20 (now 20-21) picked up an explicit
#    custom row height of 14 in the saved file.
# ------------------------------------------------------------------
$ws.Rows.Item(20).RowHeight = 14
$ws.Rows.Item(21).RowHeight = 14

# ------------------------------------------------------------------
# 6) Remove the old blank spacer row (now row 28) and make room for a
#    new validation-formula row plus keep the existing formula rows,
#    net result: a new row 30 appears, and the old formula rows 28-33
#    (now at 29-34) move down to 31-36.
# ------------------------------------------------------------------
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()

# New formula row 30
$ws.Rows.Item(30).RowHeight = 13.5
$ws.Range("B30").Formula = '=IF(AND(B19=FALSE,B18=TRUE),"This modality is not there!!!!!","ok")'

# ------------------------------------------------------------------
# 7) Move the conditional formatting that highlighted the validation
#    column from B28:B33 to the new location B30:B36.
# ------------------------------------------------------------------
$fc = $ws.Range("B28:B33").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("B30:B36"))

# ------------------------------------------------------------------
# 8) Update sheet view: remove the old scroll position / selection and
#    select C27 instead (matches the author's saved view state).
# ------------------------------------------------------------------
$ws.Range("C27").Select()
